$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.197.35"
$ws.Range("E2").Value = "  +5.04%  "
$ws.Range("D3").Value = "2.465.33"
$ws.Range("E3").Value = "  +6.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.62%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").Value = "2.464.89"
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  +5.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.08%  "
$ws.Range("D15").Value = "2.903.85"
$ws.Range("E15").Value = "  +6.46%  "
$ws.Range("D16").Value = "63.088.39"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("E17").Value = "  +6.25%  "
$ws.Range("D18").Value = "2.464.96"
$ws.Range("E18").Value = "  +6.20%  "
$ws.Range("E19").Value = "  +6.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.15%  "
$ws.Range("E21").Value = "  +5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +10.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.10%  "
$ws.Range("E30").Value = "  +13.82%  "
$ws.Range("E31").Value = "  +11.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +12.69%  "
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "369.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.66%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +12.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.21%  "
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0246"
$ws.Range("E48").Value = "  +13.23%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.48%  "
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.62%  "
